# Rewrite the feature-importance table (A2:B101) with the new feature
# names and importance scores, per the modelling/model-selection update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$featureNames = New-Object "object[,]" 100,1
$importance   = New-Object "object[,]" 100,1

$featureNames[0,0] = "aumento"
$importance[0,0] = 0.11
$featureNames[1,0] = "deferi"
$importance[1,0] = 0.04
$featureNames[2,0] = "concessão ordem"
$importance[2,0] = 0.04
$featureNames[3,0] = "inafiançabilidade"
$importance[3,0] = 0.03
$featureNames[4,0] = "medidas cautelares"
$importance[4,0] = 0.03
$featureNames[5,0] = "liminar hc"
$importance[5,0] = 0.03
$featureNames[6,0] = "xliii constituição"
$importance[6,0] = 0.02
$featureNames[7,0] = "opina deferimento"
$importance[7,0] = 0.02
$featureNames[8,0] = "concessão"
$importance[8,0] = 0.02
$featureNames[9,0] = "deferida"
$importance[9,0] = 0.02
$featureNames[10,0] = "liminar assessor"
$importance[10,0] = 0.02
$featureNames[11,0] = "deferido"
$importance[11,0] = 0.02
$featureNames[12,0] = "justiça indeferiu"
$importance[12,0] = 0.02
$featureNames[13,0] = "indeferiu liminarmente"
$importance[13,0] = 0.02
$featureNames[14,0] = "inconstitucional"
$importance[14,0] = 0.02
$featureNames[15,0] = "ficou"
$importance[15,0] = 0.02
$featureNames[16,0] = "ordem prisão"
$importance[16,0] = 0.02
$featureNames[17,0] = "opina concessão"
$importance[17,0] = 0.02
$featureNames[18,0] = "relatório http"
$importance[18,0] = 0.02
$featureNames[19,0] = "arquivado"
$importance[19,0] = 0.02
$featureNames[20,0] = "república concessão"
$importance[20,0] = 0.02
$featureNames[21,0] = "submetida"
$importance[21,0] = 0.02
$featureNames[22,0] = "precário efêmero"
$importance[22,0] = 0.01
$featureNames[23,0] = "sobrestamento"
$importance[23,0] = 0.01
$featureNames[24,0] = "xliii"
$importance[24,0] = 0.01
$featureNames[25,0] = "gabinete prestou"
$importance[25,0] = 0.01
$featureNames[26,0] = "assuntos FIANÇA"
$importance[26,0] = 0.01
$featureNames[27,0] = "implicou deferimento"
$importance[27,0] = 0.01
$featureNames[28,0] = "violência grave"
$importance[28,0] = 0.01
$featureNames[29,0] = "vedação liberdade"
$importance[29,0] = 0.01
$featureNames[30,0] = "suspender efeitos"
$importance[30,0] = 0.01
$featureNames[31,0] = "informado"
$importance[31,0] = 0.01
$featureNames[32,0] = "stj indeferiu"
$importance[32,0] = 0.01
$featureNames[33,0] = "liberdade restritiva"
$importance[33,0] = 0.01
$featureNames[34,0] = "senha primeira"
$importance[34,0] = 0.01
$featureNames[35,0] = "senha relatório"
$importance[35,0] = 0.01
$featureNames[36,0] = "outro motivo"
$importance[36,0] = 0.01
$featureNames[37,0] = "liminar espécie"
$importance[37,0] = 0.01
$featureNames[38,0] = "senha"
$importance[38,0] = 0.01
$featureNames[39,0] = "liminar suspender"
$importance[39,0] = 0.01
$featureNames[40,0] = "restritiva direitos"
$importance[40,0] = 0.01
$featureNames[41,0] = "ministro gilson"
$importance[41,0] = 0.01
$featureNames[42,0] = "mérito deste"
$importance[42,0] = 0.01
$featureNames[43,0] = "formalizado ato"
$importance[43,0] = 0.01
$featureNames[44,0] = "opinou concessão"
$importance[44,0] = 0.01
$featureNames[45,0] = "relativização"
$importance[45,0] = 0.01
$featureNames[46,0] = "fundamentos insubsistência"
$importance[46,0] = 0.01
$featureNames[47,0] = "http"
$importance[47,0] = 0.01
$featureNames[48,0] = "decisão proferida"
$importance[48,0] = 0.01
$featureNames[49,0] = "deferi pedido"
$importance[49,0] = 0.01
$featureNames[50,0] = "espécie ficou"
$importance[50,0] = 0.01
$featureNames[51,0] = "Relator_OCTAVIO GALLOTTI"
$importance[51,0] = 0.01
$featureNames[52,0] = "autoridade judiciária"
$importance[52,0] = 0.01
$featureNames[53,0] = "campo precário"
$importance[53,0] = 0.01
$featureNames[54,0] = "cautelares previstas"
$importance[54,0] = 0.01
$featureNames[55,0] = "decretada desfavor"
$importance[55,0] = 0.01
$featureNames[56,0] = "deduzida"
$importance[56,0] = 0.01
$featureNames[57,0] = "arquivado definitivo"
$importance[57,0] = 0.01
$featureNames[58,0] = "enunciado"
$importance[58,0] = 0.01
$featureNames[59,0] = "efeitos ordem"
$importance[59,0] = 0.01
$featureNames[60,0] = "enunciado súmula"
$importance[60,0] = 0.01
$featureNames[61,0] = "assuntos COMPETÊNCIA DO MP"
$importance[61,0] = 0.01
$featureNames[62,0] = "superação súmula"
$importance[62,0] = 0
$featureNames[63,0] = "requer medida"
$importance[63,0] = 0
$featureNames[64,0] = "afastamento enunciado"
$importance[64,0] = 0
$featureNames[65,0] = "resumida"
$importance[65,0] = 0
$featureNames[66,0] = "resumida prisão"
$importance[66,0] = 0
$featureNames[67,0] = "revelou contornos"
$importance[67,0] = 0
$featureNames[68,0] = "assuntos HABEAS CORPUS - LIBERATÓRIO"
$importance[68,0] = 0
$featureNames[69,0] = "Relator_ILMAR GALVÃO"
$importance[69,0] = 0
$featureNames[70,0] = "assuntos PREVISTOS NA LEGISLAÇÃO EXTRAVAGANTE"
$importance[70,0] = 0
$featureNames[71,0] = "assuntos INDEFERIMENTO"
$importance[71,0] = 0
$featureNames[72,0] = "assuntos TRANSFERÊNCIA DE PRESO"
$importance[72,0] = 0
$featureNames[73,0] = "sob código"
$importance[73,0] = 0
$featureNames[74,0] = "assuntos FURTO (ART. 155)"
$importance[74,0] = 0
$featureNames[75,0] = "análise pedido"
$importance[75,0] = 0
$featureNames[76,0] = "assuntos DESCLASSIFICAÇÃO"
$importance[76,0] = 0
$featureNames[77,0] = "suspender"
$importance[77,0] = 0
$featureNames[78,0] = "assuntos LICITAÇÕES"
$importance[78,0] = 0
$featureNames[79,0] = "ficou assim"
$importance[79,0] = 0
$featureNames[80,0] = "assim revelou"
$importance[80,0] = 0
$featureNames[81,0] = "assim resumida"
$importance[81,0] = 0
$featureNames[82,0] = "deferida assessoria"
$importance[82,0] = 0
$featureNames[83,0] = "http sob"
$importance[83,0] = 0
$featureNames[84,0] = "impetração eis"
$importance[84,0] = 0
$featureNames[85,0] = "empresas"
$importance[85,0] = 0
$featureNames[86,0] = "eis informado"
$importance[86,0] = 0
$featureNames[87,0] = "efêmero"
$importance[87,0] = 0
$featureNames[88,0] = "informado análise"
$importance[88,0] = 0
$featureNames[89,0] = "deserção"
$importance[89,0] = 0
$featureNames[90,0] = "código senha"
$importance[90,0] = 0
$featureNames[91,0] = "proferida ministro"
$importance[91,0] = 0
$featureNames[92,0] = "contornos impetração"
$importance[92,0] = 0
$featureNames[93,0] = "contornos"
$importance[93,0] = 0
$featureNames[94,0] = "cautelar pois"
$importance[94,0] = 0
$featureNames[95,0] = "aurélio decisão"
$importance[95,0] = 0
$featureNames[96,0] = "precário"
$importance[96,0] = 0
$featureNames[97,0] = "preventiva fundamentos"
$importance[97,0] = 0
$featureNames[98,0] = "processo formalizado"
$importance[98,0] = 0
$featureNames[99,0] = "assuntos DE TRÁFICO ILÍCITO E USO INDEVIDO DE DROGAS"
$importance[99,0] = 0

$ws.Range("A2:A101").Value = $featureNames
$ws.Range("B2:B101").Value = $importance

Write-Output "Updated feature importance table rows A2:B101"
